$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Union (multi-area) range of every cell whose text value changes in this
# revision (columns D "Price" and E "Volume(1h)" for the affected rows).
$changedRange = $ws.Range("D2,E2,D3,E3,D4,E4,E5,D6,E6,D7,E7,D8,E8,D9,E9,D10,E10,D11,E11,D12,D13,E13,D14,E14,D15,E15,D16,E16,D17,E17,D18,E18,D19,E19,E20,D21,E21,D22,E22,D23,E23,E24,D25,E25,E26,D27,E27,D39,E39,D40,E40,D41,E41,D42,E42,D43,E43,D44,E44,E45,D46,E46,E47,D48,E48,E49,D50,E50,D51,E51")

# Force a text number-format on every area of the range first so Excel
# keeps the assigned strings as plain text (matching the original
# t="inlineStr" cells) instead of auto-converting numeric-looking /
# percent-looking text into real numbers. NumberFormat must be applied
# area-by-area because assigning it directly on a multi-area Range only
# affects the first area.
foreach ($area in $changedRange.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "310.36"
$ws.Range("E2").Value = "1.54%"
$ws.Range("D3").Value = "35.41"
$ws.Range("E3").Value = "-2.71%"
$ws.Range("D4").Value = "5.104"
$ws.Range("E4").Value = "1.36%"
$ws.Range("E5").Value = "3.45%"
$ws.Range("D6").Value = "2.062"
$ws.Range("E6").Value = "-3.01%"
$ws.Range("D7").Value = "7.946"
$ws.Range("E7").Value = "-0.25%"
$ws.Range("D8").Value = "2.908"
$ws.Range("E8").Value = "9.30%"
$ws.Range("D9").Value = "0.9251"
$ws.Range("E9").Value = "0.23%"
$ws.Range("D10").Value = "0.1108"
$ws.Range("E10").Value = "12.83%"
$ws.Range("D11").Value = "0.1914"
$ws.Range("E11").Value = "3.10%"
$ws.Range("D12").Value = "0.09281"
$ws.Range("D13").Value = "0.03645"
$ws.Range("E13").Value = "1.31%"
$ws.Range("D14").Value = "0.09909"
$ws.Range("E14").Value = "-0.16%"
$ws.Range("D15").Value = "0.001426"
$ws.Range("E15").Value = "-0.63%"
$ws.Range("D16").Value = "0.005845"
$ws.Range("E16").Value = "4.24%"
$ws.Range("D17").Value = "3.477"
$ws.Range("E17").Value = "-0.10%"
$ws.Range("D18").Value = "4.124"
$ws.Range("E18").Value = "-0.35%"
$ws.Range("D19").Value = "0.3408"
$ws.Range("E19").Value = "-0.46%"
$ws.Range("E20").Value = "-2.04%"
$ws.Range("D21").Value = "5.098"
$ws.Range("E21").Value = "-1.47%"
$ws.Range("D22").Value = "0.2203"
$ws.Range("E22").Value = "-1.86%"
$ws.Range("D23").Value = "0.04546"
$ws.Range("E23").Value = "-0.38%"
$ws.Range("E24").Value = "-0.65%"
$ws.Range("D25").Value = "0.004814"
$ws.Range("E25").Value = "-0.34%"
$ws.Range("E26").Value = "-3.64%"
$ws.Range("D27").Value = "0.0004443"
$ws.Range("E27").Value = "-6.29%"
$ws.Range("D39").Value = "0.01973"
$ws.Range("E39").Value = "5.84%"
$ws.Range("D40").Value = "0.04875"
$ws.Range("E40").Value = "-0.49%"
$ws.Range("D41").Value = "0.007607"
$ws.Range("E41").Value = "-2.27%"
$ws.Range("D42").Value = "0.009059"
$ws.Range("E42").Value = "17.40%"
$ws.Range("D43").Value = "0.1386"
$ws.Range("E43").Value = "-0.93%"
$ws.Range("D44").Value = "0.002180"
$ws.Range("E44").Value = "-1.43%"
$ws.Range("E45").Value = "3.60%"
$ws.Range("D46").Value = "0.00006540"
$ws.Range("E46").Value = "1.98%"
$ws.Range("E47").Value = "0.07%"
$ws.Range("D48").Value = "180.86"
$ws.Range("E48").Value = "249.21%"
$ws.Range("E49").Value = "-21.03%"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "0.07%"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "0.07%"

# Restore the default "Normal" style (again area-by-area) so no stray
# style index is left behind on the cells -- the original cells carried
# no explicit style either.
foreach ($area in $changedRange.Areas) {
    $area.Style = "Normal"
}

